{"js": "// Replace the three-digit x one-digit multiplication problems throughout\n// the document's table with a new set of problems, keeping all existing\n// formatting (font, size, paragraph alignment, etc.) intact.\nconst replacements = [\n  [\"507\u00d75=2535\", \"811\u00d72=1622\"],\n  [\"674\u00d77=4718\", \"358\u00d74=1432\"],\n  [\"763\u00d79=6867\", \"134\u00d73=402\"],\n  [\"444\u00d77=3108\", \"185\u00d73=555\"],\n  [\"569\u00d74=2276\", \"468\u00d74=1872\"],\n  [\"704\u00d77=4928\", \"620\u00d79=5580\"],\n  [\"603\u00d79=5427\", \"417\u00d79=3753\"],\n  [\"542\u00d76=3252\", \"157\u00d72=314\"],\n  [\"594\u00d72=1188\", \"518\u00d73=1554\"],\n  [\"146\u00d78=1168\", \"567\u00d74=2268\"],\n  [\"933\u00d78=7464\", \"421\u00d78=3368\"],\n  [\"149\u00d74=596\", \"792\u00d75=3960\"],\n  [\"291\u00d79=2619\", \"206\u00d75=1030\"],\n  [\"402\u00d75=2010\", \"977\u00d74=3908\"],\n  [\"127\u00d79=1143\", \"675\u00d77=4725\"],\n  [\"246\u00d72=492\", \"872\u00d75=4360\"],\n  [\"679\u00d75=3395\", \"620\u00d74=2480\"],\n  [\"822\u00d73=2466\", \"434\u00d76=2604\"],\n  [\"145\u00d77=1015\", \"945\u00d74=3780\"],\n  [\"487\u00d72=974\", \"113\u00d72=226\"],\n  [\"996\u00d73=2988\", \"987\u00d73=2961\"],\n  [\"108\u00d77=756\", \"227\u00d76=1362\"],\n  [\"786\u00d78=6288\", \"792\u00d79=7128\"],\n  [\"551\u00d74=2204\", \"483\u00d72=966\"],\n  [\"966\u00d76=5796\", \"756\u00d75=3780\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems throughout\n# the document's table with a new set of problems, keeping all existing\n# formatting (font, size, paragraph alignment, etc.) intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old = \"507\u00d75=2535\"; new = \"811\u00d72=1622\"},\n    @{old = \"674\u00d77=4718\"; new = \"358\u00d74=1432\"},\n    @{old = \"763\u00d79=6867\"; new = \"134\u00d73=402\"},\n    @{old = \"444\u00d77=3108\"; new = \"185\u00d73=555\"},\n    @{old = \"569\u00d74=2276\"; new = \"468\u00d74=1872\"},\n    @{old = \"704\u00d77=4928\"; new = \"620\u00d79=5580\"},\n    @{old = \"603\u00d79=5427\"; new = \"417\u00d79=3753\"},\n    @{old = \"542\u00d76=3252\"; new = \"157\u00d72=314\"},\n    @{old = \"594\u00d72=1188\"; new = \"518\u00d73=1554\"},\n    @{old = \"146\u00d78=1168\"; new = \"567\u00d74=2268\"},\n    @{old = \"933\u00d78=7464\"; new = \"421\u00d78=3368\"},\n    @{old = \"149\u00d74=596\";  new = \"792\u00d75=3960\"},\n    @{old = \"291\u00d79=2619\"; new = \"206\u00d75=1030\"},\n    @{old = \"402\u00d75=2010\"; new = \"977\u00d74=3908\"},\n    @{old = \"127\u00d79=1143\"; new = \"675\u00d77=4725\"},\n    @{old = \"246\u00d72=492\";  new = \"872\u00d75=4360\"},\n    @{old = \"679\u00d75=3395\"; new = \"620\u00d74=2480\"},\n    @{old = \"822\u00d73=2466\"; new = \"434\u00d76=2604\"},\n    @{old = \"145\u00d77=1015\"; new = \"945\u00d74=3780\"},\n    @{old = \"487\u00d72=974\";  new = \"113\u00d72=226\"},\n    @{old = \"996\u00d73=2988\"; new = \"987\u00d73=2961\"},\n    @{old = \"108\u00d77=756\";  new = \"227\u00d76=1362\"},\n    @{old = \"786\u00d78=6288\"; new = \"792\u00d79=7128\"},\n    @{old = \"551\u00d74=2204\"; new = \"483\u00d72=966\"},\n    @{old = \"966\u00d76=5796\"; new = \"756\u00d75=3780\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $p.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
